$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.759.71'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.42%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.529.53'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.60%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.25%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.78'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("E9").Value = '  -1.27%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.13'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.96%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0810'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.46%  '

$ws.Range("E12").Value = '  -2.44%  '

$ws.Range("E13").Value = '  -3.21%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.917.78'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.63%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.547.88'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.04%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.19'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.40%  '

$ws.Range("E17").Value = '  -1.12%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.846.98'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.69%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.86'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.34%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.92'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.27%  '

$ws.Range("E21").Value = '  -0.87%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.91'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.34%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '253.17'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.90%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.08'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.88%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.71'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.20%  '

$ws.Range("E27").Value = '  +0.09%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.42'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.75%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.45'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.52%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.40'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.64%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.91'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.25%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '157.96'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.24%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.18'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.12%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.35'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.16%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.07'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.99%  '

$ws.Range("E36").Value = '  +1.83%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0783'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.38%  '

$ws.Range("E38").Value = '  -1.01%  '

$ws.Range("E39").Value = '  -1.26%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '23.25'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.73%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.31'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +13.93%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.84'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.36%  '

$ws.Range("E43").Value = '  +0.52%  '

$ws.Range("E44").Value = '  +0.29%  '

$ws.Range("E45").Value = '  -2.58%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.029.18'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.40%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.21'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.61%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.82%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '106.75'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.11%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '74.93'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.72%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.771.38'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.54%  '
